# Refresh the crypto price/volume table with the latest scraped values.
# Column D (Price) cells hold text-formatted numbers (e.g. "214.10", "0.0168")
# so we force the cell's NumberFormat to Text ("@") before assigning the
# value; otherwise Excel would silently reinterpret the text as a real
# number and drop significant trailing/leading zeros and the dotted
# thousand-separator formatting used here (e.g. "28.585.48").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.585.48"
$ws.Range("E2").Value = "  +4.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.592.66"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -1.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.10"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("E6").Value = "  +0.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -1.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.95"
$ws.Range("E8").Value = "  +8.24%  "

$ws.Range("E9").Value = "  +1.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0601"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("E11").Value = "  +2.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.815.87"
$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.590.20"
$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.78"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("E15").Value = "  +1.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.552.05"
$ws.Range("E16").Value = "  +4.93%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.97"
$ws.Range("E17").Value = "  +2.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.23"
$ws.Range("E18").Value = "  +8.74%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0710"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.42"
$ws.Range("E23").Value = "  +1.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  -0.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.46"
$ws.Range("E25").Value = "  -1.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.37"
$ws.Range("E26").Value = "  +1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.63"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.419.86"
$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  -0.29%  "

$ws.Range("E36").Value = "  -5.41%  "

$ws.Range("E37").Value = "  -1.29%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.60"
$ws.Range("E38").Value = "  +10.62%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.546"
$ws.Range("E40").Value = "  +1.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.815"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.70"
$ws.Range("E42").Value = "  -1.95%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.996"
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.84"
$ws.Range("E44").Value = "  +6.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.971"
$ws.Range("E45").Value = "  -3.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.56"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.727.33"
$ws.Range("E47").Value = "  +1.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.56"
$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0525"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0998"
$ws.Range("E50").Value = "  -4.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.93"
$ws.Range("E51").Value = "  +17.73%  "
